$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}

Replace-Text "Play Mysterious Egypt Slot for Free - Review 2021" "Play Mysterious Egypt for Free"
Replace-Text "Unique bonus feature with expanding symbols" "Excellent graphics and design"
Replace-Text "High volatility and RTP for bigger rewards" "Variety of symbols with different payout amounts"
Replace-Text "Beautiful graphics with an Ancient Egypt theme" "Unique bonus symbol that expands to cover entire reel"
Replace-Text "Flexible betting options for different strategies" "High volatility for bigger risks and rewards"
Replace-Text "Limited number of paylines may not appeal to some players" "Limited number of paylines"
Replace-Text "May not be suitable for those with a smaller bankroll" "No progressive jackpot feature"
Replace-Text "Explore the mysteries of Ancient Egypt and play Mysterious Egypt slot for free! Read our review of this high-volatility game with unique bonus features." "Read our review of Mysterious Egypt, a slot game with excellent graphics and unique bonus features. Play for free and win big!"
